$d = $word.ActiveDocument

# 1. Update the last-updated date field text
[void]$d.Content.Find.Execute("Friday, October 23, 2020", $true, $false, $false, $false, $false, $true, 1, $false, "Saturday, October 24, 2020", 2)

# 2. Merge the four runs describing built-in resource properties into a single run
$oldText = "The resource also has a number of built-in properties on the object itself. Following is the listing of those built-in properties."
[void]$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $oldText, 2)

# 3. Insert the new "Node And Socket Media Support" section right after the Resource Properties table
$tbl = $d.Tables.Item($d.Tables.Count)
$insPos = $tbl.Range.End
$rng = $d.Range($insPos, $insPos)
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/>
<w:p/>
<w:p>
  <w:pPr>
    <w:pStyle w:val="Heading2"/>
  </w:pPr>
  <w:r>
    <w:t>Node And Socket Media Support</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>Nodes and sockets both have support for multiple types of media and support for basic interaction. In general, nodes and each of a node''s sockets have support for the following elements.</w:t>
  </w:r>
</w:p>
<w:p>
  <w:r>
    <w:t>These additions allow the node and its sockets to present a self-contained card-like personality during runtime that has multiple aspects available for interactive purposes. In industry terms, this effect produces that of the Hero Card.</w:t>
  </w:r>
</w:p>
<w:p/>
<w:tbl>
  <w:tblPr>
    <w:tblStyle w:val="TableGrid"/>
    <w:tblW w:w="0" w:type="auto"/>
    <w:tblInd w:w="720" w:type="dxa"/>
    <w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>
  </w:tblPr>
  <w:tblGrid>
    <w:gridCol w:w="1201"/>
    <w:gridCol w:w="1417"/>
    <w:gridCol w:w="6012"/>
  </w:tblGrid>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D9E2F3" w:themeFill="accent1" w:themeFillTint="33"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
          <w:t>General Element</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D9E2F3" w:themeFill="accent1" w:themeFillTint="33"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
          <w:t>Property Name</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
        <w:shd w:val="clear" w:color="auto" w:fill="D9E2F3" w:themeFill="accent1" w:themeFillTint="33"/>
      </w:tcPr>
      <w:p>
        <w:pPr>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
        </w:pPr>
        <w:r>
          <w:rPr>
            <w:b/>
            <w:bCs/>
          </w:rPr>
          <w:t>Description</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Audio</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>MediaAudio</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>On nodes, this is a URI referring to an audio to play when the node is loaded. On sockets, the URI is played when the option is selected and before the connection is followed to the next node.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:lastRenderedPageBreak/>
          <w:t>Image</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>MediaImage</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>A URI referring to an image to display in the card. base64 data URI is allowed.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Link</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>MediaLink</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t xml:space="preserve">A URL the user can click to view more information about the current context. Standard </w:t>
        </w:r>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>href</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
        <w:r>
          <w:t xml:space="preserve"> link syntax.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
  <w:tr>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>Video</w:t>
        </w:r>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:proofErr w:type="spellStart"/>
        <w:r>
          <w:t>MediaVideo</w:t>
        </w:r>
        <w:proofErr w:type="spellEnd"/>
      </w:p>
    </w:tc>
    <w:tc>
      <w:tcPr>
        <w:tcW w:w="0" w:type="auto"/>
      </w:tcPr>
      <w:p>
        <w:r>
          <w:t>On nodes, this is a URI referring to a video to play when the node is loaded. On sockets, the URI is played when the option is selected, and before the connection is followed to the next node.</w:t>
        </w:r>
      </w:p>
    </w:tc>
  </w:tr>
</w:tbl>
<w:p/>
<w:p>
  <w:r>
    <w:t>Notice that in this version, if audio and video URIs are both specified, the video URI will take precedence, and the audio will only be utilized on non-video systems that provide audio support, such as a telephone system or radio, etc.</w:t>
  </w:r>
</w:p>
<w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$rng.InsertXML($xml)

Write-Host "done"
